{"js": "// \"added intro dialogs to level 2\"\n//\n// 1) Paragraph \"So the weather ... later.\" loses its gramStart/gramEnd\n//    proofing-error markers and the leading \"So\" run is merged with the\n//    rest of the sentence into a single run.\n// 2) \"Though beetles ... seem to prefer eating plants.\" gets corrected to\n//    \"...seems to prefer eating plants.\".\n// 3) A whole new block of dialog paragraphs (climate oceanic / microclimate /\n//    mushroom) is appended after \"These insects thrive in most climates...\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraphs by their text instead of hard-coded indexes\n// so the script is resilient to any minor paragraph-count drift.\nconst items = paragraphs.items;\n\nconst weatherParaIndex = items.findIndex((p) =>\n  p.text.indexOf(\"the weather we will be experiencing here over a few days\") !== -1\n);\nconst beetleParaIndex = items.findIndex((p) =>\n  p.text.indexOf(\"Though beetles almost eat anything\") !== -1\n);\nconst lastInsectParaIndex = items.findIndex((p) =>\n  p.text.indexOf(\"These insects thrive in most climates\") !== -1\n);\n\nif (weatherParaIndex === -1 || beetleParaIndex === -1 || lastInsectParaIndex === -1) {\n  throw new Error(\"Could not find the expected anchor paragraphs.\");\n}\n\n// 1) Rebuild the \"So the weather ...\" paragraph as a single clean run, which\n//    also drops the w:proofErr gramStart/gramEnd markers around \"So\".\nconst weatherPara = items[weatherParaIndex];\nconst nextAfterWeather = items[weatherParaIndex + 1];\nnextAfterWeather.insertParagraph(\n  \"So the weather we will be experiencing here over a few days will be vastly different several months later.\",\n  Word.InsertLocation.before\n);\nweatherPara.delete();\n\nawait context.sync();\n\n// 2) Fix \"seem\" -> \"seems\" in the beetle paragraph.\nconst beetlePara = items[beetleParaIndex];\nconst seemResults = beetlePara.search(\"seem\", { matchCase: true });\nseemResults.load(\"items\");\nawait context.sync();\nseemResults.items[0].insertText(\"s\", Word.InsertLocation.after);\nawait context.sync();\n\n// 3) Append the new intro-dialog paragraphs after \"These insects thrive...\".\nconst lastInsectPara = items[lastInsectParaIndex];\n\nconst newParagraphTexts = [\n  \"\",\n  \"(climate oceanic)\",\n  \"Ah, the oceanic climate! Where it\\u2019s nice and cool throughout the year!\",\n  \"This is due to the region being close to the ocean that regulates the temperature.\",\n  \"However, the winds carrying the nice ocean breeze can cause a lot of days to be dull and dreary.\",\n  \"Nevertheless, this moody atmosphere is a perfect place to grow our plants!\",\n  \"\",\n  \"(microclimate)\",\n  \"Looks like we landed in a peculiar place where the weather is a little bit cooler than usual.\",\n  \"This is known as a microclimate, where the local atmospheric condition can differ from the surrounding areas.\",\n  \"In our case, we are below the slope of a hill that obscures most of the sunlight, and precipitations linger around longer.\",\n  \"\",\n  \"(mushroom)\",\n  \"Uh oh, a mushroom started to grow near our plants!\",\n  \"Since there\\u2019s a lot of moisture in the region, the fungi that grow these mushrooms from beneath are able to absorb a lot of nutrients.\",\n  \"Their spores appear to be harmful to all our frogs and plants, but one: the iron frog!\",\n  \"Just as the iron frogs can rid us of the moles, they, too, can rid us of these mushrooms.\"\n];\n\nlet anchorPara = lastInsectPara;\nfor (const text of newParagraphTexts) {\n  anchorPara = anchorPara.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# \"added intro dialogs to level 2\"\n#\n# 1) Paragraph \"So the weather ... later.\" loses its gramStart/gramEnd\n#    proofing-error markers and the leading \"So\" run is merged with the\n#    rest of the sentence into a single run.\n# 2) \"Though beetles ... seem to prefer eating plants.\" gets corrected to\n#    \"...seems to prefer eating plants.\".\n# 3) A whole new block of dialog paragraphs (climate oceanic / microclimate /\n#    mushroom) is appended after \"These insects thrive in most climates...\".\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraphs by text instead of a hard-coded index so the\n# script is resilient to any minor paragraph-count drift.\n$weatherIndex = 0\n$beetleIndex = 0\n$lastInsectIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*the weather we will be experiencing here over a few days*\") {\n        $weatherIndex = $i\n    }\n    if ($t -like \"*Though beetles almost eat anything*\") {\n        $beetleIndex = $i\n    }\n    if ($t -like \"*These insects thrive in most climates*\") {\n        $lastInsectIndex = $i\n    }\n}\n\n# 1) Rebuild the \"So the weather ...\" paragraph as a single clean run, which\n#    also drops the w:proofErr gramStart/gramEnd markers around \"So\".\n$weatherPara = $d.Paragraphs.Item($weatherIndex)\n$nextPara = $d.Paragraphs.Item($weatherIndex + 1)\n$insertRange = $nextPara.Range\n$insertRange.Collapse(1)  # wdCollapseStart\n$insertRange.InsertBefore(\"So the weather we will be experiencing here over a few days will be vastly different several months later.`r\")\n$weatherPara.Range.Delete()\n\n# 2) Fix \"seem\" -> \"seems\" in the beetle paragraph.\n$beetlePara = $d.Paragraphs.Item($beetleIndex)\n$findRange = $beetlePara.Range\n$find = $findRange.Find\n$find.Text = \"seem\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n$findRange.Collapse(0)  # wdCollapseEnd\n$findRange.InsertAfter(\"s\")\n\n# 3) Append the new intro-dialog paragraphs after \"These insects thrive...\".\n$newParagraphTexts = @(\n    \"\",\n    \"(climate oceanic)\",\n    \"Ah, the oceanic climate! Where it\u2019s nice and cool throughout the year!\",\n    \"This is due to the region being close to the ocean that regulates the temperature.\",\n    \"However, the winds carrying the nice ocean breeze can cause a lot of days to be dull and dreary.\",\n    \"Nevertheless, this moody atmosphere is a perfect place to grow our plants!\",\n    \"\",\n    \"(microclimate)\",\n    \"Looks like we landed in a peculiar place where the weather is a little bit cooler than usual.\",\n    \"This is known as a microclimate, where the local atmospheric condition can differ from the surrounding areas.\",\n    \"In our case, we are below the slope of a hill that obscures most of the sunlight, and precipitations linger around longer.\",\n    \"\",\n    \"(mushroom)\",\n    \"Uh oh, a mushroom started to grow near our plants!\",\n    \"Since there\u2019s a lot of moisture in the region, the fungi that grow these mushrooms from beneath are able to absorb a lot of nutrients.\",\n    \"Their spores appear to be harmful to all our frogs and plants, but one: the iron frog!\",\n    \"Just as the iron frogs can rid us of the moles, they, too, can rid us of these mushrooms.\"\n)\n\n$lastPara = $d.Paragraphs.Item($lastInsectIndex)\nforeach ($t in $newParagraphTexts) {\n    $r = $lastPara.Range\n    $r.Collapse(0)  # wdCollapseEnd\n    $r.InsertParagraphAfter()\n    $lastPara = $d.Paragraphs.Last\n    if ($t -ne \"\") {\n        $lastPara.Range.InsertAfter($t)\n    }\n}\n"}
